# Append new row 89 to the Optical_Power worksheet, matching the data
# collected by the automated map update (2025-07-24 06:58:13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A, B and D hold numeric/date-looking values but are stored as
# plain text throughout the sheet (e.g. "-529", "7/23/2025", "1").
# Force a text number format before writing so the engine doesn't
# auto-convert them to a number/date, then restore the default style so
# no extra formatting is left behind on the cell.
$ws.Range("A89").NumberFormat = "@"
$ws.Range("A89").Value = "-529"
$ws.Range("A89").Style = "Normal"

$ws.Range("B89").NumberFormat = "@"
$ws.Range("B89").Value = "7/23/2025"
$ws.Range("B89").Style = "Normal"

$ws.Range("C89").Value = "Libertad 820"

$ws.Range("D89").NumberFormat = "@"
$ws.Range("D89").Value = "1"
$ws.Range("D89").Style = "Normal"

# E89 is blank (no OT value for this case) in the source data.
$ws.Range("E89").NumberFormat = "@"
$ws.Range("E89").Value = ""
$ws.Range("E89").Style = "Normal"

$ws.Range("F89").Value = "Optical Power"
$ws.Range("G89").Value = "Pendiente"
$ws.Range("H89").Value = "Colocar columna hablar con Pablo si hay dudas"
$ws.Range("I89").Value = 1
$ws.Range("J89").Value = "Cambio"
$ws.Range("K89").Value = "Sin equipos"
$ws.Range("L89").Value = "Pasante"
$ws.Range("M89").Value = -58.384097
$ws.Range("N89").Value = -34.598913
$ws.Range("O89").Value = "Recoleta"
$ws.Range("P89").Value = "Capital Sur"
